$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.435.05"
$ws.Range("E2").Value = "  -2.67%  "

$ws.Range("D3").Value = "3.174.54"
$ws.Range("E3").Value = "  -4.36%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "'569.93"
$ws.Range("E5").Value = "  -3.07%  "

$ws.Range("D6").Value = "'168.69"
$ws.Range("E6").Value = "  -7.90%  "

$ws.Range("E7").Value = "  -6.09%  "

$ws.Range("D9").Value = "3.174.34"
$ws.Range("E9").Value = "  -4.20%  "

$ws.Range("D10").Value = "'0.121"
$ws.Range("E10").Value = "  -4.12%  "

$ws.Range("D11").Value = "'6.79"
$ws.Range("E11").Value = "  +0.01%  "

$ws.Range("D12").Value = "'0.386"
$ws.Range("E12").Value = "  -3.84%  "

$ws.Range("D13").Value = "3.721.91"
$ws.Range("E13").Value = "  -4.39%  "

$ws.Range("E14").Value = "  -2.25%  "

$ws.Range("D15").Value = "64.459.50"
$ws.Range("E15").Value = "  -2.67%  "

$ws.Range("D16").Value = "'25.35"
$ws.Range("E16").Value = "  -3.37%  "

$ws.Range("E17").Value = "  -2.54%  "

$ws.Range("D18").Value = "3.164.88"
$ws.Range("E18").Value = "  -3.53%  "

$ws.Range("D19").Value = "'417.87"
$ws.Range("E19").Value = "  -2.58%  "

$ws.Range("D20").Value = "'5.37"
$ws.Range("E20").Value = "  -2.86%  "

$ws.Range("D21").Value = "'12.85"
$ws.Range("E21").Value = "  -2.62%  "

$ws.Range("D22").Value = "'7.08"
$ws.Range("E22").Value = "  -4.57%  "

$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  +0.18%  "

$ws.Range("D24").Value = "'69.99"
$ws.Range("E24").Value = "  -2.63%  "

$ws.Range("D25").Value = "'0.205"
$ws.Range("E25").Value = "  -0.53%  "

$ws.Range("D26").Value = "'0.487"
$ws.Range("E26").Value = "  -5.82%  "

$ws.Range("E27").Value = "  -7.41%  "

$ws.Range("D28").Value = "'8.86"
$ws.Range("E28").Value = "  -0.71%  "

$ws.Range("E29").Value = "  +0.17%  "

$ws.Range("E30").Value = "  -5.91%  "

$ws.Range("D31").Value = "'21.65"
$ws.Range("E31").Value = "  -3.24%  "

$ws.Range("E32").Value = "  -0.07%  "

$ws.Range("D33").Value = "'5.04"
$ws.Range("E33").Value = "  -2.78%  "

$ws.Range("D34").Value = "'6.33"
$ws.Range("E34").Value = "  -3.83%  "

$ws.Range("D35").Value = "'1.14"
$ws.Range("E35").Value = "  -4.22%  "

$ws.Range("D36").Value = "'157.47"
$ws.Range("E36").Value = "  -1.38%  "

$ws.Range("E37").Value = "  -5.65%  "

$ws.Range("D38").Value = "2.728.59"
$ws.Range("E38").Value = "  -5.46%  "

$ws.Range("E39").Value = "  -5.73%  "

$ws.Range("D40").Value = "'24.32"
$ws.Range("E40").Value = "  -8.32%  "

$ws.Range("D41").Value = "'4.18"
$ws.Range("E41").Value = "  -3.31%  "

$ws.Range("D42").Value = "'39.14"
$ws.Range("E42").Value = "  -2.28%  "

$ws.Range("D43").Value = "'0.710"
$ws.Range("E43").Value = "  -7.19%  "

$ws.Range("D44").Value = "'0.0622"
$ws.Range("E44").Value = "  -6.51%  "

$ws.Range("D45").Value = "'5.62"
$ws.Range("E45").Value = "  -4.90%  "

$ws.Range("D46").Value = "'0.0262"
$ws.Range("E46").Value = "  -3.81%  "

$ws.Range("D47").Value = "'293.71"
$ws.Range("E47").Value = "  -6.91%  "

$ws.Range("D48").Value = "'21.59"
$ws.Range("E48").Value = "  -6.99%  "

$ws.Range("B49").Value = "dogwifhat"
$ws.Range("C49").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D49").Value = "'2.01"
$ws.Range("E49").Value = "  -12.58%  "

$ws.Range("B50").Value = "FirstDigitalUSD"
$ws.Range("C50").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D50").Value = "'1.00"
$ws.Range("E50").Value = "  +0.06%  "

$ws.Range("E51").Value = "  -6.25%  "
